$d = $word.ActiveDocument
Write-Output $d.Styles.Count
